# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 4
    4  = 5
    5  = 5
    6  = 8
    7  = 7
    8  = 8
    9  = 3
    10 = 6
    11 = 10
    12 = 8
    13 = 4
    14 = 12
    15 = 8
    16 = 8
    17 = 10
    18 = 8
    19 = 8
    20 = 7
    21 = 4
    22 = 3
    23 = 8
    24 = 7
    25 = 9
    26 = 8
    27 = 9
    28 = 5
    29 = 4
    30 = 5
    31 = 7
    32 = 6
    33 = 4
    34 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
